$p = $ppt.ActivePresentation

# Slide 1: title "First slide" -- collapse the word-by-word runs into a
# single run. Delete()+InsertAfter() (rather than a plain Text= assignment,
# which would be a silent no-op since the concatenated text is unchanged)
# forces the run-merge and keeps the paragraph's pPr intact.
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Delete() | Out-Null
$tr1.InsertAfter("First slide") | Out-Null

# Slide 2 (blank slide): its notes page body placeholder had the same text
# split word-by-word across many runs; collapse it to a single run.
$s2 = $p.Slides.Item(2)
$tr2 = $s2.NotesPage.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Some notes here: this first slide should use the Blank template"

# Slide 3: title "Third slide" -- same run-collapsing as slide 1.
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Delete() | Out-Null
$tr3.InsertAfter("Third slide") | Out-Null
